$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 4872.8
$ws.Range("I6").Value = 407.8
$ws.Range("K6").Value = 1223.4
$ws.Range("M6").Value = -1111.4
$ws.Range("H8").Value = 12405614
$ws.Range("I8").Value = 12405614
$ws.Range("K8").Value = 37216842
$ws.Range("M8").Value = -37216703
$ws.Range("H12").Value = 183.16667
$ws.Range("I12").Value = 243
$ws.Range("J12").Value = 123.333336
$ws.Range("K12").Value = 243
$ws.Range("L12").Value = 123.333336
$ws.Range("M12").Value = -73
$ws.Range("N12").Value = -463.333336
$ws.Range("H64").Value = 6233.222
$ws.Range("I64").Value = 3500
$ws.Range("K64").Value = 3500
$ws.Range("M64").Value = -3252
$ws.Range("H67").Value = 6233.222
$ws.Range("I67").Value = 3500
$ws.Range("K67").Value = 3500
$ws.Range("M67").Value = -2642
$ws.Range("H74").Value = 3179.3157
$ws.Range("I74").Value = 2386.2144
$ws.Range("J74").Value = 5400
$ws.Range("K74").Value = 2386.2144
$ws.Range("L74").Value = 5400
$ws.Range("M74").Value = -1450.2144
$ws.Range("N74").Value = -7272
$ws.Range("H77").Value = 3179.3157
$ws.Range("I77").Value = 2386.2144
$ws.Range("J77").Value = 5400
$ws.Range("K77").Value = 11931.072
$ws.Range("L77").Value = 27000
$ws.Range("M77").Value = -7251.072
$ws.Range("N77").Value = -36360
$ws.Range("H112").Value = 1715.0667
$ws.Range("J112").Value = 1817.4615
$ws.Range("L112").Value = 5452.3845
$ws.Range("N112").Value = -7668.3845
$ws.Range("H125").Value = 1405.625
$ws.Range("I125").Value = 1299.5
$ws.Range("J125").Value = 1420.7858
$ws.Range("K125").Value = 11695.5
$ws.Range("L125").Value = 12787.0722
$ws.Range("M125").Value = -9235.5
$ws.Range("N125").Value = -17707.0722
$ws.Range("H127").Value = 453250
$ws.Range("I127").Value = 453250
$ws.Range("K127").Value = 1359750
$ws.Range("M127").Value = -1354790
$ws.Range("H131").Value = 1109.3
$ws.Range("I131").Value = 1200
$ws.Range("K131").Value = 3600
$ws.Range("M131").Value = 1440
$ws.Range("H133").Value = 52172
$ws.Range("J133").Value = 56296
$ws.Range("L133").Value = 56296
$ws.Range("N133").Value = -66416
$ws.Range("H138").Value = 1647.0741
$ws.Range("J138").Value = 3009.875
$ws.Range("L138").Value = 9029.625
$ws.Range("N138").Value = -19309.625

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 37037.35
$ws.Range("I32").Value = 21962.146
$ws.Range("K32").Value = 21962.146
$ws.Range("M32").Value = -21675.146
$ws.Range("H45").Value = 14686.885
$ws.Range("I45").Value = 13831.556
$ws.Range("K45").Value = 13831.556
$ws.Range("M45").Value = -13454.556
$ws.Range("H74").Value = 1134.9736
$ws.Range("I74").Value = 798.0357
$ws.Range("J74").Value = 2078.4
$ws.Range("K74").Value = 798.0357
$ws.Range("L74").Value = 2078.4
$ws.Range("M74").Value = 75.96429999999998
$ws.Range("N74").Value = -3826.4
$ws.Range("H77").Value = 1134.9736
$ws.Range("I77").Value = 798.0357
$ws.Range("J77").Value = 2078.4
$ws.Range("K77").Value = 3990.1785
$ws.Range("L77").Value = 10392
$ws.Range("M77").Value = 377.8215
$ws.Range("N77").Value = -19128
$ws.Range("H132").Value = 2025.6
$ws.Range("I132").Value = 1455.7826
$ws.Range("J132").Value = 3117.75
$ws.Range("K132").Value = 4367.3478
$ws.Range("L132").Value = 9353.25
$ws.Range("M132").Value = -1837.3478
$ws.Range("N132").Value = -14413.25

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 6952908.5
$ws.Range("I20").Value = 15878779
$ws.Range("K20").Value = 15878779
$ws.Range("M20").Value = -15878532
$ws.Range("H94").Value = 1115.2142
$ws.Range("I94").Value = 482.66666
$ws.Range("J94").Value = 2253.8
$ws.Range("K94").Value = 482.66666
$ws.Range("L94").Value = 2253.8
$ws.Range("M94").Value = -31.66665999999998
$ws.Range("N94").Value = -3155.8

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 413.0435
$ws.Range("I7").Value = 298
$ws.Range("K7").Value = 298
$ws.Range("M7").Value = -185
$ws.Range("H31").Value = 4262.278
$ws.Range("I31").Value = 1616.1666
$ws.Range("J31").Value = 5585.3335
$ws.Range("K31").Value = 1616.1666
$ws.Range("L31").Value = 5585.3335
$ws.Range("M31").Value = -1321.1666
$ws.Range("N31").Value = -6175.3335
$ws.Range("H34").Value = 4262.278
$ws.Range("I34").Value = 1616.1666
$ws.Range("J34").Value = 5585.3335
$ws.Range("K34").Value = 1616.1666
$ws.Range("L34").Value = 5585.3335
$ws.Range("M34").Value = -1414.1666
$ws.Range("N34").Value = -5989.3335
$ws.Range("H132").Value = 1010.5
$ws.Range("I132").Value = 1010.5
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 3031.5
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -501.5
$ws.Range("N132").ClearContents()

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H20").Value = 466.33334
$ws.Range("J20").Value = 466.33334
$ws.Range("L20").Value = 1399.00002
$ws.Range("N20").Value = -1853.00002

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 18185840
$ws.Range("I70").Value = 22226028
$ws.Range("K70").Value = 22226028
$ws.Range("M70").Value = -22225758
$ws.Range("H73").Value = 18185840
$ws.Range("I73").Value = 22226028
$ws.Range("K73").Value = 22226028
$ws.Range("M73").Value = -22225092
$ws.Range("H132").Value = 1553.3334
$ws.Range("I132").Value = 1578
$ws.Range("J132").Value = 1504
$ws.Range("K132").Value = 4734
$ws.Range("L132").Value = 4512
$ws.Range("M132").Value = -2204
$ws.Range("N132").Value = -9572

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1366.8
$ws.Range("I16").Value = 1366.8
$ws.Range("K16").Value = 1366.8
$ws.Range("M16").Value = -1196.8
$ws.Range("H40").Value = 5400.5884
$ws.Range("I40").Value = 4901.143
$ws.Range("K40").Value = 4901.143
$ws.Range("M40").Value = -4765.143
$ws.Range("H46").Value = 1744.6923
$ws.Range("I46").Value = 1590.6666
$ws.Range("K46").Value = 1590.6666
$ws.Range("M46").Value = -1402.6666
$ws.Range("H132").Value = 2596.2424
$ws.Range("I132").Value = 2307.64
$ws.Range("J132").Value = 3498.125
$ws.Range("K132").Value = 6922.92
$ws.Range("L132").Value = 10494.375
$ws.Range("M132").Value = -4392.92
$ws.Range("N132").Value = -15554.375
